# LAP22G33-115 US301 UCD, MD and SSD
# Applies the self-assessment / rubric updates described by the commit.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet: Group and Self Assessment
# ---------------------------------------------------------------------------
$wsGroup = $wb.Worksheets.Item("Group and Self Assessment")

$wsGroup.Range("B4").Value = 33

# Row 10 (student 1201237)
$wsGroup.Range("D10").Value = 5
$wsGroup.Range("E10").Value = 5
$wsGroup.Range("F10").Value = 4
$wsGroup.Range("G10").Value = 4

# Row 11 (student 1201274)
$wsGroup.Range("D11").Value = 5
$wsGroup.Range("E11").Value = 5
$wsGroup.Range("G11").Value = 4

# Row 12 (student 1190772)
$wsGroup.Range("D12").Value = 5
$wsGroup.Range("E12").Value = 5
$wsGroup.Range("G12").Value = 4

# Row 13 (student 1200546)
$wsGroup.Range("F13").Value = 4
$wsGroup.Range("G13").Value = 3

# ---------------------------------------------------------------------------
# Sheet: User Stories
# ---------------------------------------------------------------------------
$wsUS = $wb.Worksheets.Item("User Stories")

$wsUS.Range("A6").Value = "US107"
$wsUS.Range("B6").Value = 1201237
$wsUS.Range("C6").Value = 4

$wsUS.Range("A7").Value = "US201"

$wsUS.Range("A8").Value = "US202"
$wsUS.Range("B8").Value = 1201274

$wsUS.Range("A9").Value = "US203"
$wsUS.Range("B9").Value = 1201237
$wsUS.Range("C9").Value = 4

$wsUS.Range("A10").Value = "US204"
$wsUS.Range("D10").ClearContents()

$wsUS.Range("A11").Value = "US205"
$wsUS.Range("B11").Value = 1200546
$wsUS.Range("C11").Value = 1

$wsUS.Range("A12").Value = "US206"
$wsUS.Range("C12").Value = 1

$wsUS.Range("A13").Value = "US207"

$wsUS.Range("A14").Value = "US208"
$wsUS.Range("B14").Value = 1190772

$wsUS.Range("A15").Value = "US209"

$wsUS.Range("A16").Value = "US210"
$wsUS.Range("B16").Value = 1201274
$wsUS.Range("C16").Value = 2

$wsUS.Range("A17").Value = "US111"
$wsUS.Range("B17").Value = 1200546
$wsUS.Range("C17").Value = 4

# ---------------------------------------------------------------------------
# Sheet: Code Quality
# ---------------------------------------------------------------------------
$wsCQ = $wb.Worksheets.Item("Code Quality")

$wsCQ.Range("C4").Value = 83.9
$wsCQ.Range("C5").Value = 73.9
$wsCQ.Range("D5").Value = 65
$wsCQ.Range("C6").Value = 0.9

# ---------------------------------------------------------------------------
# Sheet: Project Development
# ---------------------------------------------------------------------------
$wsPD = $wb.Worksheets.Item("Project Development")

$wsPD.Range("F5").Value = 3
$wsPD.Range("D6").Value = 3
$wsPD.Range("F6").Value = 3
$wsPD.Range("D7").Value = 2
$wsPD.Range("F7").Value = 1
$wsPD.Range("F8").Value = 4

# ---------------------------------------------------------------------------
# Sheet: Project Management
# ---------------------------------------------------------------------------
$wsPM = $wb.Worksheets.Item("Project Management")

$wsPM.Range("D5").Value = 2
$wsPM.Range("E5").Value = 2
$wsPM.Range("F5").Value = 2

$wsPM.Range("C7").Value = 4
$wsPM.Range("D7").Value = 4
$wsPM.Range("E7").Value = 3
$wsPM.Range("F7").Value = 1

$wsPM.Range("D9").Value = 5
$wsPM.Range("E9").Value = 3
$wsPM.Range("F9").Value = 2

$wsPM.Range("C11").Value = 2
$wsPM.Range("D11").Value = 2
$wsPM.Range("E11").Value = 2
$wsPM.Range("F11").Value = 2

$wsPM.Range("E13").Value = 5
$wsPM.Range("F13").Value = 1

# ---------------------------------------------------------------------------
# Sheet views / selections (cosmetic, mirrors the authored commit) and
# establishing "Group and Self Assessment" as the final active sheet/tab.
# ---------------------------------------------------------------------------
$wsUS.Activate()
$wsUS.Range("C12").Select()

$wsCQ.Activate()
$wsCQ.Range("D11").Select()

$wsPD.Activate()
$wsPD.Range("F7").Select()

$wsPM.Activate()
$wsPM.Range("F13").Select()

$wsGroup.Activate()
$wsGroup.Range("H13").Select()
